# LOQ4055.docx restructuring edit.
#
# The paragraph/run *structure* (styles, bold/italic runs, line breaks)
# stays the same; only the text carried by certain runs has been reshuffled
# between paragraphs/slots. So instead of re-arranging paragraphs we move
# the text itself, slot by slot.
#
# Several of the moves form cycles (A's old text becomes B's new text while
# B's old text becomes A's new text, etc.), so everything is staged through
# unique placeholder tokens first and only resolved to final text in a
# second pass -- this guarantees no step can accidentally match text a
# previous step just wrote.

$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $null = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# ---------------------------------------------------------------------
# Whole-paragraph slots (each of these paragraphs carries exactly one
# run, so replacing the paragraph Range.Text is safe and keeps the
# paragraph's own pPr/rPr formatting, e.g. the italic runs).
# ---------------------------------------------------------------------

# 1) Objetivos (PT) -> Programa resumido (PT) short text
$d.Paragraphs(6).Range.Text = "@@PH_6@@"
# 2) Objetivos (EN, italic) -> Programa resumido (EN, italic) short text
$d.Paragraphs(7).Range.Text = "@@PH_7@@"
# 3) Docente(s) bullet -> old Objetivos (PT) text
$d.Paragraphs(9).Range.Text = "@@PH_9@@"
# 4) Programa resumido (PT) -> Programa (PT) long text
$d.Paragraphs(11).Range.Text = "@@PH_11@@"
# 5) Programa resumido (EN, italic) -> old Objetivos (EN) text
$d.Paragraphs(12).Range.Text = "@@PH_12@@"
# 6) Programa (PT) -> Método short text
$d.Paragraphs(14).Range.Text = "@@PH_14@@"
# 7) Old Bibliografia body paragraph -> Docente bullet content
$d.Paragraphs(19).Range.Text = "@@PH_19@@"

# ---------------------------------------------------------------------
# Runs inside the multi-run "Avaliação" paragraph (#17): only the
# non-bold value runs move; the bold label runs ("Método: ", "Critério: ",
# "Norma de recuperação: ") stay untouched.
# ---------------------------------------------------------------------
Replace-All "Serão oferecidas aulas expositivas." "@@PH_M@@"
Replace-All "Serão aplicadas duas provas escritas. Trabalhos em sala de aula, seminários e relatórios, poderão, a critério do docente, ser considerados como parte da nota da prova escrita." "@@PH_C@@"
Replace-All "Será realizada uma prova escrita envolvendo o conteúdo do semestre todo." "@@PH_N@@"

# ---------------------------------------------------------------------
# Phase 2: resolve every placeholder to its real, final text.
# ---------------------------------------------------------------------

Replace-All "@@PH_6@@" "Métodos de separação e obtenção dos elementos, extração mineral. Hidrogênio. Metais alcalinos. Metais alcalino terrosos. Alumínio. Metais de transição. Compostos de coordenação. Halogênios."
Replace-All "@@PH_7@@" "Methods for separating and obtaining elements, mineral extraction. Hydrogen. Alkaline metals. Alkaline earth metals. Aluminum. Transition metals. Coordination compounds. Halogens."
Replace-All "@@PH_9@@" "Fornecer aos alunos conceitos fundamentos para a compreensão da Química Inorgânica, de forma a capacitá-lo a descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de caráter inorgânico com interesse industrial."
Replace-All "@@PH_11@@" "Métodos de separação e obtenção dos elementos, extração mineral. Propriedades, obtenção e aplicações dos seguintes elementos/grupos e seus compostos: Hidrogênio; Metais alcalinos (indústria cloro-álcali; processo Solvay); Metais alcalino terrosos; Alumínio (processo Bayer); Metais de transição; Compostos de coordenação e Halogênios."
Replace-All "@@PH_12@@" "Provide students with fundamental concepts for understanding Inorganic Chemistry, in order to enable them to describe and interpret the properties of elements and their compounds, especially those of an inorganic nature with industrial interest"
Replace-All "@@PH_14@@" "Serão oferecidas aulas expositivas."
Replace-All "@@PH_19@@" "5840705 - Maria Lúcia Caetano Pinto da Silva"

Replace-All "@@PH_M@@" "Serão aplicadas duas provas escritas. Trabalhos em sala de aula, seminários e relatórios, poderão, a critério do docente, ser considerados como parte da nota da prova escrita."
Replace-All "@@PH_C@@" "Será realizada uma prova escrita envolvendo o conteúdo do semestre todo."

# "Norma de recuperação: " run now absorbs the whole bibliography list,
# preserving the manual line breaks (double) originally between entries.
# A backtick-v is PowerShell's escape for a vertical-tab char, which Word
# treats as a manual line break (<w:br/>) when assigned into Range.Text.
$bib = "WELLER, Mark; OVERTON, Tina; ROURKE, Jonathan; et al. Química inorgânica. Porto Alegre, Bookman, 6ª Ed, 2017. E-book. " + "`v`v" + "CHANG, Raymond. Química geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010." + "`v`v" + "BROWN, T.L. ET al. Química a ciência central. 9.ed. São Paulo: Pearson Prentice Hall, 2005-2007." + "`v`v" + "BRADY, J ; HUMISTON, G.E. Química geral. Rio de Janeiro: Ed. Livros Técnicos Científicos, 1981." + "`v`v" + "LEE, J. D., tradução Química Inorgânica não tão concisa da 5ª edição inglesa. Editora Edgard Blucher Ltda. SP-2001." + "`v`v" + "SHRIVER, D. e ATKINS, P. Química Inorgânica tradução da 4ª edição. Editora Bookman, Porto Alegre-RS, 2008." + "`v`v" + "QUAGLIANO, J.V; VALLARINO, L.M. Química - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3ª ed., 1973."
Replace-All "@@PH_N@@" $bib
